# Auto-generated edit script: weekly refresh of Berenjena price rows
# (rows shift down by two weeks of history; two new latest-week rows
# appended; dataset grows from 44 to 46 data rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns, in order: Row, A Mercado ID, B Mercado, C Region, D Fecha(serial),
# E Codreg, F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
# K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg, Q Kg o Unidades, R Clasificacion
$data = @(
  @(16, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44804, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 50, 12000, 12000, 12000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 200, 60, "Hortaliza"),
  @(17, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44804, 16, 100112001, "Berenjena", "Sin especificar", "Segunda", 30, 12000, 12000, 12000, "`$/caja 90 unidades", "Región de Arica y Parinacota", 133, 90, "Hortaliza"),
  @(18, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44610, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 11000, 12000, 11500, "`$/caja 60 unidades", "Región Metropolitana", 192, 60, "Hortaliza"),
  @(19, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44698, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 10000, 10000, 10000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 167, 60, "Hortaliza"),
  @(20, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44218, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 65, 9000, 10000, 9615, "`$/caja 60 unidades", "Región del Maule", 160, 60, "Hortaliza"),
  @(21, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44792, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 12000, 13000, 12500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 208, 60, "Hortaliza"),
  @(22, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44204, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 45, 9500, 10000, 9722, "`$/caja 60 unidades", "Región del Maule", 162, 60, "Hortaliza"),
  @(23, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44755, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 11000, 12000, 11500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 192, 60, "Hortaliza"),
  @(24, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44160, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 90, 7500, 8000, 7667, "`$/caja 60 unidades", "Región de Arica y Parinacota", 128, 60, "Hortaliza"),
  @(25, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44271, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 55, 9000, 9500, 9227, "`$/caja 60 unidades", "Región del Maule", 154, 60, "Hortaliza"),
  @(26, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44224, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 80, 8500, 9000, 8719, "`$/caja 60 unidades", "Región del Maule", 145, 60, "Hortaliza"),
  @(27, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44790, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 12000, 13000, 12500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 208, 60, "Hortaliza"),
  @(28, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44615, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 100, 11000, 12000, 11500, "`$/caja 60 unidades", "Región Metropolitana", 192, 60, "Hortaliza"),
  @(29, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44784, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 80, 12000, 13000, 12500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 208, 60, "Hortaliza"),
  @(30, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44776, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 11000, 12000, 11500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 192, 60, "Hortaliza"),
  @(31, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44594, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 80, 12000, 13000, 12500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 208, 60, "Hortaliza"),
  @(32, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44216, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 55, 9500, 10000, 9773, "`$/caja 60 unidades", "Región del Maule", 163, 60, "Hortaliza"),
  @(33, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44671, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 160, 6000, 7000, 6500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 108, 60, "Hortaliza"),
  @(34, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44692, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 10000, 10000, 10000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 167, 60, "Hortaliza"),
  @(35, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44259, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 70, 9000, 9500, 9214, "`$/caja 60 unidades", "Región del Maule", 154, 60, "Hortaliza"),
  @(36, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44627, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 9000, 9500, 9250, "`$/caja 60 unidades", "Región Metropolitana", 154, 60, "Hortaliza"),
  @(37, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44264, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 43, 8500, 9000, 8709, "`$/caja 60 unidades", "Región del Maule", 145, 60, "Hortaliza"),
  @(38, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44764, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 12000, 13000, 12500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 208, 60, "Hortaliza"),
  @(39, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44761, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 12000, 13000, 12500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 208, 60, "Hortaliza"),
  @(40, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44798, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 40, 12000, 12000, 12000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 200, 60, "Hortaliza"),
  @(41, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44600, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 12000, 13000, 12500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 208, 60, "Hortaliza"),
  @(42, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44763, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 12000, 13000, 12500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 208, 60, "Hortaliza"),
  @(43, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44699, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 10000, 10000, 10000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 167, 60, "Hortaliza"),
  @(44, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44202, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 50, 8000, 9000, 8400, "`$/caja 60 unidades", "Región del Maule", 140, 60, "Hortaliza"),
  @(45, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44159, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 35, 7500, 8000, 7714, "`$/caja 60 unidades", "Región de Arica y Parinacota", 129, 60, "Hortaliza"),
  @(46, 7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44782, 16, 100112001, "Berenjena", "Sin especificar", "Primera", 60, 12000, 13000, 12500, "`$/caja 60 unidades", "Región de Arica y Parinacota", 208, 60, "Hortaliza"),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
    $ws.Cells.Item($r, 16).Value = $row[16]
    $ws.Cells.Item($r, 17).Value = $row[17]
    $ws.Cells.Item($r, 18).Value = $row[18]
}
